# Update "想去人数" (want-to-go count) values in column F
# on the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$sheetA = $wb.Worksheets.Item("展览")
$sheetB = $wb.Worksheets.Item("全部类型")

# Updates for "展览" sheet (F column, keyed by row number)
$updatesA = @{
    2  = 88
    3  = 11912
    4  = 20
    6  = 355
    8  = 11817
    9  = 491
    10 = 1174
    11 = 102
    12 = 65
    13 = 1781
    14 = 5855
    16 = 3539
    17 = 190
    18 = 22
}

foreach ($row in $updatesA.Keys) {
    $sheetA.Range("F$row").Value = $updatesA[$row]
}

# Updates for "全部类型" sheet (F column, keyed by row number)
$updatesB = @{
    3  = 88
    5  = 11912
    6  = 20
    9  = 355
    11 = 11817
    12 = 491
    13 = 1174
    14 = 102
    15 = 65
    16 = 1781
    18 = 5855
    20 = 3539
    21 = 190
    22 = 22
}

foreach ($row in $updatesB.Keys) {
    $sheetB.Range("F$row").Value = $updatesB[$row]
}
